$p = $ppt.ActivePresentation

# --- Slide 1: was empty, now gets a Title + Content placeholder ("Group 1 project") ---
$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)
$s1.CustomLayout = $s2.CustomLayout
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Group 1 project "

# --- Slide 2: "General results" + empty content placeholder  -->
#     becomes "ggseqlogo" title, with the picture that used to live on slide 3 ---
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "ggseqlogo"
$s2.Shapes.Item(2).Delete()

$s3 = $p.Slides.Item(3)
$pic = $s3.Shapes.Item(2)
$pic.Copy()
$s2.Shapes.Paste() | Out-Null

# --- Build the brand-new "Other questions" slide by duplicating slide 1 (it already
#     has the plain Title+Content placeholder shapes we want). This changes slide
#     indices from here on, so slides are re-fetched by Item(...) as needed. ---
$s1 = $p.Slides.Item(1)
$new = $s1.Duplicate()
$new.Shapes.Item(1).TextFrame.TextRange.Text = "Other questions "
$body = $new.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Force to ignore errors in source I Rmd "
$body.InsertAfter("`rShiny " + [char]0x2013 + " insert blank option ? ")
$body.InsertAfter("`rDo you want anything specific according to the github page ? ")
$body.InsertAfter("`rWhat to put in a packages ? ")

# Move the new slide to the end of the deck.
$new.MoveTo($p.Slides.Count)

# --- Slide 3 ("Questions.... What to do :/") is removed: its title text was folded
#     into slide 2 above, and its picture was copied onto slide 2 above. The old
#     slide 4 ("Modelling ?? :/") shifts up and becomes the new slide 3. ---
$p.Slides.Item(3).Delete()
